$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94
$ws.Cells.Item(94, 1).Value = 1337640217444924928
$ws.Cells.Item(94, 2).Value = 'Kita kasih waktu sampai malam tahun baru, deh~
Cek Instagram @kelaskitadotcom sekarang juga!
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #ikancupang #giveaway #parangkencana https://t.co/Cxu1EhOwe9'
$ws.Cells.Item(94, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(94, 4).Value = 'Sat Dec 12 06:07:34 +0000 2020'

# Row 95
$ws.Cells.Item(95, 1).Value = 1337563060463161088
$ws.Cells.Item(95, 2).Value = 'Punya topik atau pertanyaan menarik untuk dibahas? Ayo sharing bareng kami di sesi lia s. bookclub! 📚 Tayang live di Instagram setiap Sabtu jam 11.00 WIB, jangan sampai ketinggalan ya 😃
#liasidikbranding #belajaronline #books #growth #branding101 https://t.co/kXMtmE4GJ2'
$ws.Cells.Item(95, 3).Value = 'liasidik'
$ws.Cells.Item(95, 4).Value = 'Sat Dec 12 01:00:58 +0000 2020'

# Row 96
$ws.Cells.Item(96, 1).Value = 1337370303165465088
$ws.Cells.Item(96, 2).Value = '(..con) m dari https://t.co/yzZy8M21Lb.
#belajardaring
#belajaronline
#kknupi2020
#dirumahaja https://t.co/ZfXXd4UULF'
$ws.Cells.Item(96, 3).Value = 'Krsmynt1'
$ws.Cells.Item(96, 4).Value = 'Fri Dec 11 12:15:02 +0000 2020'

# Row 97
$ws.Cells.Item(97, 1).Value = 1337335960791446016
$ws.Cells.Item(97, 2).Value = 'Ini hari, Kelaskita mau ngasih tau siapa-siapa aja ''Pang Jago yang berhasil mendapatkan hadiah keren dari @kelaskitadotcom.
Cek Instagram @kelaskitadotcom sekarang!
#kelaskita #carabarubelajarseru #HaloPangJago #belajardirumah #belajaronline #dirumahaja #ikancupang #giveaway https://t.co/hcHpLGsXrS'
$ws.Cells.Item(97, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(97, 4).Value = 'Fri Dec 11 09:58:34 +0000 2020'

# Row 98
$ws.Cells.Item(98, 1).Value = 1337265588230450944
$ws.Cells.Item(98, 2).Value = 'Pastikan anak anda memiliki cangkir yang penuh! 
#Vyneapple #funlearning #keluarga #studygram #belajar #belajaronline #parenthood #parenting #parentingtips #aplikasiandroid #aplikasiios #pendidikan #anakpintar #edukasianak https://t.co/zcwxAF3s2V'
$ws.Cells.Item(98, 3).Value = 'vyneapple'
$ws.Cells.Item(98, 4).Value = 'Fri Dec 11 05:18:56 +0000 2020'

$ws.Range("K94").Select()
